$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.86000000000045
$ws.Range("H2").Value = [double]"1.197005956469171e-16"
$ws.Range("K2").Value = 45.27693583471693
$ws.Range("L2").Value = "[42.31162049110234, 48.24225117833152]"
$ws.Range("O2").Value = 1.616395018964117
$ws.Range("P2").Value = "[1.5409213215805773, 1.6918687163476571]"
$ws.Range("S2").Value = 50.01543416943114
$ws.Range("T2").Value = "[48.03244034766444, 51.998427991197836]"
$ws.Range("W2").Value = 18.46458458458492
$ws.Range("X2").Value = 18.16596596596629
$ws.Range("Y2").Value = 18.76320320320355

# Row 3 updates
$ws.Range("E3").Value = 23.00000000000016
$ws.Range("H3").Value = [double]"1.197005956469171e-16"
$ws.Range("I3").Value = 0.93894575122854
$ws.Range("K3").Value = 48.14815891776909
$ws.Range("L3").Value = "[43.7178600958609, 52.578457739677276]"
$ws.Range("S3").Value = 50.81400155247482
$ws.Range("T3").Value = "[47.73664857537094, 53.8913545295787]"
$ws.Range("W3").Value = 5.755755755755796
$ws.Range("X3").Value = 5.387387387387424
$ws.Range("Y3").Value = 6.124124124124169
